$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string but must remain plain text
# (they were, and must stay, text cells - e.g. "6.49", "20.46" ...).
# We temporarily force a text number-format so Excel does not auto-convert the
# assigned string into a real number, then restore each cell's original style.
$forceTextCells = @("D5", "D6", "D8", "D10", "D12", "D13", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D43", "D44", "D46", "D47", "D48", "D51")

$savedStyles = @{}
foreach ($addr in $forceTextCells) {
    $cell = $ws.Range($addr)
    $savedStyles[$addr] = $cell.Style
    $cell.NumberFormat = "@"
}

# Apply all the updated values from the source diff.
$ws.Range("D2").Value = "58.764.00"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.587.09"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "519.72"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "139.72"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "2.600.33"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "0.331"
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").Value = "0.134"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "3.044.38"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "58.816.25"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "20.46"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.576.65"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "339.07"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "4.32"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "10.20"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "6.52"
$ws.Range("E22").Value = "  +5.73%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "66.22"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "7.08"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").Value = "5.95"
$ws.Range("E31").Value = "  -4.72%  "
$ws.Range("D32").Value = "18.77"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").Value = "1.56"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "148.87"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "3.98"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "36.27"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("D38").Value = "0.834"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").Value = "0.818"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "274.87"
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").Value = "10.75"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "0.586"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "0.0521"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "18.57"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "1.982.43"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "4.47"
$ws.Range("E51").Value = "  -0.97%  "

# Restore original styles on the cells we temporarily reformatted.
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = $savedStyles[$addr]
}
